$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header relabel: "Node ID" -> "Item" ---
$ws.Range("A1").Value = "Item"

# --- "Batman Robin" -> "Batman" ---
$ws.Range("K2").Value = "Batman"

# --- Row 9 field-name headers relabeled to a more generic set ---
$ws.Range("B9").Value = "Hello"
$ws.Range("C9").Value = "Yeah"
$ws.Range("D9").Value = "Its"
$ws.Range("E9").Value = "Been"
$ws.Range("F9").Value = "A_While"
$ws.Range("G9").Value = "Not"
$ws.Range("H9").Value = "Much"
$ws.Range("I9").Value = "How"
$ws.Range("J9").Value = "About"
$ws.Range("K9").Value = "You"

# --- Numeric values bumped ---
$ws.Range("G1").Value = 111000
$ws.Range("H1").Value = 111000

# --- View state: scroll back to the top-left and select row 34 ---
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows(34).Select()

# --- Column H width nudged down a hair ---
$ws.Columns("H").ColumnWidth = 12.05
